# Weekly price-sheet update: insert a new "Cilantro" record for
# Terminal La Palmera de La Serena as row 138 (the most recent week,
# 2022-08-25 == serial 44798), pushing the existing historical rows
# 138-163 down to 139-164.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 138; Excel shifts rows 138:163 down to 139:164
# and the sheet's used range grows to A1:R164 automatically.
$ws.Rows(138).Insert()

$ws.Cells.Item(138, 1).Value  = 8
$ws.Cells.Item(138, 2).Value  = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(138, 3).Value  = 'Coquimbo'
$ws.Cells.Item(138, 4).Value  = 44798
$ws.Cells.Item(138, 5).Value  = 4
$ws.Cells.Item(138, 6).Value  = 100112040
$ws.Cells.Item(138, 7).Value  = 'Cilantro'
$ws.Cells.Item(138, 8).Value  = 'Sin especificar'
$ws.Cells.Item(138, 9).Value  = 'Primera'
$ws.Cells.Item(138, 10).Value = 2400
$ws.Cells.Item(138, 11).Value = 2000
$ws.Cells.Item(138, 12).Value = 2500
$ws.Cells.Item(138, 13).Value = 2250
$ws.Cells.Item(138, 14).Value = '$/atado 1 a 1,5 kilos'
$ws.Cells.Item(138, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(138, 16).Value = 1500
$ws.Cells.Item(138, 17).Value = 1.5
$ws.Cells.Item(138, 18).Value = 'Hortaliza'
